# Auto-generated edit script applying updated market/profit data
# to specific rows across multiple worksheets, per commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 3802.5
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 3802.5
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 11407.5
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -14901.5

$ws.Range("H132").Value = 6384.02
$ws.Range("I132").Value = 5768.1177
$ws.Range("J132").Value = 7692.8125
$ws.Range("K132").Value = 17304.3531
$ws.Range("L132").Value = 23078.4375
$ws.Range("M132").Value = -14774.3531
$ws.Range("N132").Value = -28138.4375

$ws.Range("H137").Value = 1600.85
$ws.Range("I137").Value = 1923.7826
$ws.Range("J137").Value = 1163.9412
$ws.Range("K137").Value = 5771.3478
$ws.Range("L137").Value = 3491.8236
$ws.Range("M137").Value = -3221.3478
$ws.Range("N137").Value = -8591.8236

$ws.Range("H138").Value = 1463.7413
$ws.Range("I138").Value = 1193.4736
$ws.Range("J138").Value = 1977.25
$ws.Range("K138").Value = 3580.4208
$ws.Range("L138").Value = 5931.75
$ws.Range("M138").Value = 1559.5792
$ws.Range("N138").Value = -16211.75

$ws.Range("H141").Value = 8370.817999999999
$ws.Range("I141").Value = 3509.9333
$ws.Range("J141").Value = 18787
$ws.Range("K141").Value = 10529.7999
$ws.Range("L141").Value = 56361
$ws.Range("M141").Value = -5349.7999
$ws.Range("N141").Value = -66721

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6618.202
$ws.Range("I32").Value = 5797.857
$ws.Range("J32").Value = 20400
$ws.Range("K32").Value = 5797.857
$ws.Range("L32").Value = 20400
$ws.Range("M32").Value = -5510.857
$ws.Range("N32").Value = -20974

$ws.Range("H45").Value = 1819.5385
$ws.Range("I45").Value = 1578
$ws.Range("J45").Value = 1926.8889
$ws.Range("K45").Value = 1578
$ws.Range("L45").Value = 1926.8889
$ws.Range("M45").Value = -1201
$ws.Range("N45").Value = -2680.8889

$ws.Range("H74").Value = 2165.3845
$ws.Range("I74").Value = 2112.5
$ws.Range("J74").Value = 2250
$ws.Range("K74").Value = 2112.5
$ws.Range("L74").Value = 2250
$ws.Range("M74").Value = -1238.5
$ws.Range("N74").Value = -3998

$ws.Range("H77").Value = 2165.3845
$ws.Range("I77").Value = 2112.5
$ws.Range("J77").Value = 2250
$ws.Range("K77").Value = 10562.5
$ws.Range("L77").Value = 11250
$ws.Range("M77").Value = -6194.5
$ws.Range("N77").Value = -19986

$ws.Range("H97").Value = 961.4286
$ws.Range("I97").Value = 961.4286
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 961.4286
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -465.4286
$ws.Range("N97").ClearContents()

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H132").Value = 886237.4
$ws.Range("I132").Value = 2509484.2
$ws.Range("J132").Value = 5046.2285
$ws.Range("K132").Value = 7528452.600000001
$ws.Range("L132").Value = 15138.6855
$ws.Range("M132").Value = -7525922.600000001
$ws.Range("N132").Value = -20198.6855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1393.0588
$ws.Range("I94").Value = 1217.6666
$ws.Range("K94").Value = 1217.6666
$ws.Range("M94").Value = -766.6666

$ws.Range("H133").Value = 25632
$ws.Range("J133").Value = 25632
$ws.Range("L133").Value = 25632
$ws.Range("N133").Value = -35752

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3000.5625
$ws.Range("I16").Value = 3234.25
$ws.Range("J16").Value = 2299.5
$ws.Range("K16").Value = 3234.25
$ws.Range("L16").Value = 2299.5
$ws.Range("M16").Value = -2947.25
$ws.Range("N16").Value = -2873.5

$ws.Range("H31").Value = 5052857.5
$ws.Range("I31").Value = 1877.079
$ws.Range("J31").Value = 11907760
$ws.Range("K31").Value = 1877.079
$ws.Range("L31").Value = 11907760
$ws.Range("M31").Value = -1582.079
$ws.Range("N31").Value = -11908350

$ws.Range("H34").Value = 5052857.5
$ws.Range("I34").Value = 1877.079
$ws.Range("J34").Value = 11907760
$ws.Range("K34").Value = 1877.079
$ws.Range("L34").Value = 11907760
$ws.Range("M34").Value = -1675.079
$ws.Range("N34").Value = -11908164

$ws.Range("H58").Value = 7325
$ws.Range("I58").Value = 3376.375
$ws.Range("J58").Value = 10834.889
$ws.Range("K58").Value = 3376.375
$ws.Range("L58").Value = 10834.889
$ws.Range("M58").Value = -3173.375
$ws.Range("N58").Value = -11240.889

$ws.Range("H68").Value = 9923
$ws.Range("I68").Value = 9923
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 9923
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -9174
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 9923
$ws.Range("I71").Value = 9923
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 29769
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -26025
$ws.Range("N71").ClearContents()

$ws.Range("H105").Value = 1432.1666
$ws.Range("I105").Value = 1432.1666
$ws.Range("K105").Value = 1432.1666
$ws.Range("M105").Value = 314.8334

$ws.Range("H107").Value = 1589.0667
$ws.Range("I107").Value = 285.2857
$ws.Range("K107").Value = 285.2857
$ws.Range("M107").Value = 1634.7143

$ws.Range("H113").Value = 3000.5625
$ws.Range("I113").Value = 3234.25
$ws.Range("J113").Value = 2299.5
$ws.Range("K113").Value = 3234.25
$ws.Range("L113").Value = 2299.5
$ws.Range("M113").Value = -1064.25
$ws.Range("N113").Value = -6639.5

$ws.Range("H132").Value = 3286.7585
$ws.Range("I132").Value = 3380.4
$ws.Range("K132").Value = 10141.2
$ws.Range("M132").Value = -7611.200000000001

$ws.Range("H136").Value = 7325
$ws.Range("I136").Value = 3376.375
$ws.Range("J136").Value = 10834.889
$ws.Range("K136").Value = 10129.125
$ws.Range("L136").Value = 32504.667
$ws.Range("M136").Value = -7579.125
$ws.Range("N136").Value = -37604.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 14708468
$ws.Range("I126").Value = 25001702
$ws.Range("J126").Value = 3848.8572
$ws.Range("K126").Value = 75005106
$ws.Range("L126").Value = 11546.5716
$ws.Range("M126").Value = -75002636
$ws.Range("N126").Value = -16486.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1227.1428
$ws.Range("I82").Value = 816.875
$ws.Range("J82").Value = 2540
$ws.Range("K82").Value = 816.875
$ws.Range("L82").Value = 2540
$ws.Range("M82").Value = -455.875
$ws.Range("N82").Value = -3262

$ws.Range("H85").Value = 1227.1428
$ws.Range("I85").Value = 816.875
$ws.Range("J85").Value = 2540
$ws.Range("K85").Value = 816.875
$ws.Range("L85").Value = 2540
$ws.Range("M85").Value = 431.125
$ws.Range("N85").Value = -5036

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 29389
$ws.Range("J93").Value = 29389
$ws.Range("L93").Value = 29389
$ws.Range("N93").Value = -34381

$ws.Range("H132").Value = 2654.8518
$ws.Range("I132").Value = 2534.818
$ws.Range("J132").Value = 2737.375
$ws.Range("K132").Value = 7604.454000000001
$ws.Range("L132").Value = 8212.125
$ws.Range("M132").Value = -5074.454000000001
$ws.Range("N132").Value = -13272.125
